$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix "Marking" row (per-question marks)
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Fix "Total" row (overall totals)
$ws.Range("B12").Value = 56
$ws.Range("C12").Value = -2
$ws.Range("E12").Value = "54 / 112"
